$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5822.8
$ws.Range("J64").Value = 5485
$ws.Range("L64").Value = 5485
$ws.Range("N64").Value = -5981

$ws.Range("H67").Value = 5822.8
$ws.Range("J67").Value = 5485
$ws.Range("L67").Value = 5485
$ws.Range("N67").Value = -7201

$ws.Range("H98").Value = 1687.4546
$ws.Range("I98").Value = 1766.15
$ws.Range("J98").Value = 900.5
$ws.Range("K98").Value = 1766.15
$ws.Range("L98").Value = 900.5
$ws.Range("M98").Value = -268.1500000000001
$ws.Range("N98").Value = -3896.5

$ws.Range("H105").Value = 26223
$ws.Range("J105").Value = 26223
$ws.Range("L105").Value = 26223
$ws.Range("N105").Value = -33211

$ws.Range("H116").Value = 207422.53
$ws.Range("I116").Value = 379542.38
$ws.Range("J116").Value = 10714.143
$ws.Range("K116").Value = 379542.38
$ws.Range("L116").Value = 10714.143
$ws.Range("M116").Value = -376100.38
$ws.Range("N116").Value = -17598.143

$ws.Range("H122").Value = 1687.4546
$ws.Range("I122").Value = 1766.15
$ws.Range("J122").Value = 900.5
$ws.Range("K122").Value = 5298.450000000001
$ws.Range("L122").Value = 2701.5
$ws.Range("M122").Value = -2848.450000000001
$ws.Range("N122").Value = -7601.5

$ws.Range("H132").Value = 4081.1587
$ws.Range("I132").Value = 2593.2075
$ws.Range("J132").Value = 11967.3
$ws.Range("K132").Value = 7779.622499999999
$ws.Range("L132").Value = 35901.89999999999
$ws.Range("M132").Value = -5249.622499999999
$ws.Range("N132").Value = -40961.89999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 252242
$ws.Range("J45").Value = 4226.9414
$ws.Range("L45").Value = 4226.9414
$ws.Range("N45").Value = -4980.9414

$ws.Range("H60").Value = 20249.5
$ws.Range("I60").Value = 20249.5
$ws.Range("K60").Value = 20249.5
$ws.Range("M60").Value = -19516.5

$ws.Range("H61").Value = 3258.25
$ws.Range("I61").Value = 1877.7778
$ws.Range("K61").Value = 1877.7778
$ws.Range("M61").Value = -1665.7778

$ws.Range("H88").Value = 1761.9231
$ws.Range("I88").Value = 660.2
$ws.Range("K88").Value = 660.2
$ws.Range("M88").Value = -254.2

$ws.Range("H91").Value = 1761.9231
$ws.Range("I91").Value = 660.2
$ws.Range("K91").Value = 660.2
$ws.Range("M91").Value = 743.8

$ws.Range("H122").Value = 1674.2
$ws.Range("I122").Value = 1716.8235
$ws.Range("J122").Value = 1432.6666
$ws.Range("K122").Value = 5150.470499999999
$ws.Range("L122").Value = 4297.9998
$ws.Range("M122").Value = -2700.470499999999
$ws.Range("N122").Value = -9197.9998

$ws.Range("H132").Value = 38631.785
$ws.Range("I132").Value = 41395.23
$ws.Range("J132").Value = 2707
$ws.Range("K132").Value = 124185.69
$ws.Range("L132").Value = 8121
$ws.Range("M132").Value = -121655.69
$ws.Range("N132").Value = -13181

$ws.Range("H136").Value = 3258.25
$ws.Range("I136").Value = 1877.7778
$ws.Range("K136").Value = 5633.3334
$ws.Range("M136").Value = -3083.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35302.332
$ws.Range("I35").Value = 35302.332
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 35302.332
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -34992.332
$ws.Range("N35").ClearContents()

$ws.Range("H94").Value = 1280.3125
$ws.Range("I94").Value = 957.1667
$ws.Range("J94").Value = 2249.75
$ws.Range("K94").Value = 957.1667
$ws.Range("L94").Value = 2249.75
$ws.Range("M94").Value = -506.1667
$ws.Range("N94").Value = -3151.75

$ws.Range("H134").Value = 2275.647
$ws.Range("I134").Value = 2044.2128
$ws.Range("J134").Value = 4995
$ws.Range("K134").Value = 6132.6384
$ws.Range("L134").Value = 14985
$ws.Range("M134").Value = -3597.6384
$ws.Range("N134").Value = -20055

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 22309.846
$ws.Range("I58").Value = 27949.625
$ws.Range("J58").Value = 3510.5833
$ws.Range("K58").Value = 27949.625
$ws.Range("L58").Value = 3510.5833
$ws.Range("M58").Value = -27746.625
$ws.Range("N58").Value = -3916.5833

$ws.Range("H99").Value = 3781.7896
$ws.Range("I99").Value = 3150.4614
$ws.Range("J99").Value = 5149.6665
$ws.Range("K99").Value = 3150.4614
$ws.Range("L99").Value = 5149.6665
$ws.Range("M99").Value = -1652.4614
$ws.Range("N99").Value = -8145.6665

$ws.Range("H126").Value = 3781.7896
$ws.Range("I126").Value = 3150.4614
$ws.Range("J126").Value = 5149.6665
$ws.Range("K126").Value = 9451.3842
$ws.Range("L126").Value = 15448.9995
$ws.Range("M126").Value = -6981.3842
$ws.Range("N126").Value = -20388.9995

$ws.Range("H132").Value = 2946.8518
$ws.Range("I132").Value = 2814.6956
$ws.Range("J132").Value = 3706.75
$ws.Range("K132").Value = 8444.086800000001
$ws.Range("L132").Value = 11120.25
$ws.Range("M132").Value = -5914.086800000001
$ws.Range("N132").Value = -16180.25

$ws.Range("H136").Value = 22309.846
$ws.Range("I136").Value = 27949.625
$ws.Range("J136").Value = 3510.5833
$ws.Range("K136").Value = 83848.875
$ws.Range("L136").Value = 10531.7499
$ws.Range("M136").Value = -81298.875
$ws.Range("N136").Value = -15631.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 1070
$ws.Range("I24").Value = 712.5
$ws.Range("K24").Value = 2137.5
$ws.Range("M24").Value = -1907.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1778.0358
$ws.Range("I102").Value = 1309.909
$ws.Range("J102").Value = 3494.5
$ws.Range("K102").Value = 1309.909
$ws.Range("L102").Value = 3494.5
$ws.Range("M102").Value = 312.0909999999999
$ws.Range("N102").Value = -6738.5

$ws.Range("H132").Value = 54100.5
$ws.Range("I132").Value = 65577.56
$ws.Range("K132").Value = 196732.68
$ws.Range("M132").Value = -194202.68

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 48933.258
$ws.Range("I132").Value = 52312.04
$ws.Range("J132").Value = 6698.5
$ws.Range("K132").Value = 156936.12
$ws.Range("L132").Value = 20095.5
$ws.Range("M132").Value = -154406.12
$ws.Range("N132").Value = -25155.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 535092.0600000001
$ws.Range("I81").Value = 1242.1428
$ws.Range("J81").Value = 846504.5
$ws.Range("K81").Value = 2484.2856
$ws.Range("L81").Value = 1693009
$ws.Range("M81").Value = -1423.2856
$ws.Range("N81").Value = -1695131

$ws.Range("H84").Value = 535092.0600000001
$ws.Range("I84").Value = 1242.1428
$ws.Range("J84").Value = 846504.5
$ws.Range("K84").Value = 12421.428
$ws.Range("L84").Value = 8465045
$ws.Range("M84").Value = -7117.428
$ws.Range("N84").Value = -8475653

$ws.Range("H107").Value = 448.73685
$ws.Range("I107").Value = 251.8125
$ws.Range("J107").Value = 1499
$ws.Range("K107").Value = 755.4375
$ws.Range("L107").Value = 4497
$ws.Range("M107").Value = 1164.5625
$ws.Range("N107").Value = -8337

$ws.Range("H126").Value = 145492.86
$ws.Range("I126").Value = 145492.86
$ws.Range("K126").Value = 436478.58
$ws.Range("M126").Value = -434008.58

$ws.Range("H132").Value = 135721.89
$ws.Range("I132").Value = 140187.12
$ws.Range("K132").Value = 420561.36
$ws.Range("M132").Value = -418031.36
